$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("Q3:Q7").Value = "Liam"
$ws.Range("Q8:Q12").Value = "Brailey"

$ws.Range("Q8").Select()
$ws.Range("Q8:Q12").Select()

$ws.Application.ActiveWindow.ScrollColumn = 6
